$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.649.19'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.977.89'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.21'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.81%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.380'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0791'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.05%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.28'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.844'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('D14').Value = '2.266.06'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.79'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.43'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.91%  '
$ws.Range('D17').Value = '1.976.70'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '36.583.60'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.91'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.147'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.20'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.40'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.42'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +19.58%  '
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.81'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0617'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.52'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.29'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.40%  '
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.41'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -11.01%  '
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.98'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').Value = '1.367.68'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.22'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.21'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.83'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.48%  '
$ws.Range('D51').Value = '2.159.34'
$ws.Range('E51').Value = '  +0.55%  '
